$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.818.67"
$ws.Range("E2").Value = "  +0.54%  "

# Row 3
$ws.Range("D3").Value = "1.645.57"

# Row 4
$ws.Range("E4").Value = "  +0.46%  "

# Row 5
$ws.Range("E5").Value = "  +0.78%  "

# Row 6
$ws.Range("E6").Value = "  -0.56%  "

# Row 7
$ws.Range("E7").Value = "  +0.48%  "

# Row 8
$ws.Range("E8").Value = "  -0.20%  "

# Row 9
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.21"
$ws.Range("E10").Value = "  -0.35%  "

# Row 11
$ws.Range("E11").Value = "  +0.43%  "

# Row 12
$ws.Range("D12").Value = "1.661.43"
$ws.Range("E12").Value = "  +1.16%  "

# Row 13
$ws.Range("E13").Value = "  -0.67%  "

# Row 14
$ws.Range("E14").Value = "  -0.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.72"
$ws.Range("E15").Value = "  -0.64%  "

# Row 16
$ws.Range("D16").Value = "26.812.25"
$ws.Range("E16").Value = "  +0.48%  "

# Row 17
$ws.Range("E17").Value = "  -1.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.73"
$ws.Range("E18").Value = "  -0.80%  "

# Row 19
$ws.Range("E19").Value = "  +0.47%  "

# Row 20
$ws.Range("E20").Value = "  +0.79%  "

# Row 21
$ws.Range("E21").Value = "  +8.47%  "

# Row 22
$ws.Range("E22").Value = "  -0.23%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.32"
$ws.Range("E23").Value = "  -1.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "146.12"
$ws.Range("E24").Value = "  +0.16%  "

# Row 25
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("E26").Value = "  -1.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.19"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.69"

# Row 29
$ws.Range("E29").Value = "  -1.23%  "

# Row 30
$ws.Range("E30").Value = "  +0.69%  "

# Row 31
$ws.Range("E31").Value = "  -0.74%  "

# Row 32
$ws.Range("E32").Value = "  -1.43%  "

# Row 33
$ws.Range("D33").Value = "1.284.96"
$ws.Range("E33").Value = "  +0.70%  "

# Row 34
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("E35").Value = "  +1.59%  "

# Row 36
$ws.Range("E36").Value = "  -0.98%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.537"
$ws.Range("E37").Value = "  +0.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.821"
$ws.Range("E38").Value = "  -1.34%  "

# Row 39
$ws.Range("E39").Value = "  +0.59%  "

# Row 40
$ws.Range("E40").Value = "  -1.15%  "

# Row 41
$ws.Range("E41").Value = "  -0.51%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.32"
$ws.Range("E42").Value = "  -2.44%  "

# Row 43
$ws.Range("D43").Value = "1.785.51"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.55"
$ws.Range("E44").Value = "  +3.12%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.87"
$ws.Range("E45").Value = "  +0.61%  "

# Row 46
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  +1.01%  "

# Row 48
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0518"
$ws.Range("E48").Value = "  +0.54%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.67"
$ws.Range("E49").Value = "  -2.04%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0971"
$ws.Range("E50").Value = "  +0.23%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  +0.25%  "
